$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new row of data (row 34)
# Order of assignment matters for shared-string table ordering: Term (B) first, then Organ (A), then Variant (C), then Link (D)
$ws.Range("B34").Value = "Adenoma (Primary hyperparathyroidism) "
$ws.Range("A34").Value = "Parathyroid "
$ws.Range("C34").Value = "Clip 1 B-mode + Color "
$ws.Range("D34").Value = "https://youtu.be/S45odD2wQOQ"

# Add hyperlink on D34 and apply hyperlink style
$ws.Hyperlinks.Add($ws.Range("D34"), "https://youtu.be/S45odD2wQOQ") | Out-Null
$ws.Range("D34").Style = "Collegamento ipertestuale"

# Update dimension / view happens automatically; adjust selection to mimic final state
$ws.Range("D38").Select()
